{"js": "// Replace the twenty-five \"NNN\u00d7N=\" multiplication prompts in the body's\n// table cells with their new values. Each \"from\" string is unique in the\n// original document (verified against the source OOXML), and no \"from\"\n// string is ever the same as an earlier \"to\" string in this list, so a\n// straightforward ordered sequence of exact, case-sensitive searches +\n// replacements is safe even though one new value (\"202\u00d73=\") duplicates an\n// original value used elsewhere in the table.\nconst replacements = [\n  [\"514\u00d76=\", \"985\u00d74=\"],\n  [\"846\u00d75=\", \"299\u00d74=\"],\n  [\"169\u00d78=\", \"369\u00d77=\"],\n  [\"123\u00d76=\", \"937\u00d74=\"],\n  [\"202\u00d73=\", \"831\u00d77=\"],\n  [\"544\u00d73=\", \"515\u00d75=\"],\n  [\"975\u00d77=\", \"402\u00d79=\"],\n  [\"411\u00d72=\", \"321\u00d77=\"],\n  [\"474\u00d78=\", \"965\u00d73=\"],\n  [\"930\u00d73=\", \"868\u00d78=\"],\n  [\"939\u00d74=\", \"916\u00d77=\"],\n  [\"114\u00d78=\", \"718\u00d77=\"],\n  [\"144\u00d78=\", \"110\u00d73=\"],\n  [\"570\u00d76=\", \"437\u00d79=\"],\n  [\"325\u00d74=\", \"818\u00d79=\"],\n  [\"788\u00d78=\", \"721\u00d73=\"],\n  [\"895\u00d76=\", \"310\u00d79=\"],\n  [\"130\u00d73=\", \"606\u00d76=\"],\n  [\"647\u00d75=\", \"466\u00d76=\"],\n  [\"340\u00d79=\", \"202\u00d73=\"],\n  [\"944\u00d78=\", \"575\u00d76=\"],\n  [\"250\u00d76=\", \"555\u00d76=\"],\n  [\"757\u00d72=\", \"668\u00d79=\"],\n  [\"556\u00d77=\", \"679\u00d72=\"],\n  [\"167\u00d79=\", \"576\u00d76=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [from, to] of replacements) {\n  const found = body.search(from, { matchCase: true, matchWholeWord: false });\n  found.load(\"text\");\n  await context.sync();\n\n  if (found.items.length === 0) {\n    throw new Error(`Could not find text to replace: ${from}`);\n  }\n\n  // Each source string is unique in the document, but replace defensively\n  // over every hit in case the runtime ever reports more than one.\n  for (let i = 0; i < found.items.length; i++) {\n    found.items[i].insertText(to, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the twenty-five \"NNN\u00d7N=\" multiplication prompts in the body's\n# table cells with their new values. Each \"from\" string is unique in the\n# original document (verified against the source OOXML), and no \"from\"\n# string is ever the same as an earlier \"to\" string in this list, so a\n# straightforward ordered sequence of exact, case-sensitive Find/Replace\n# calls is safe even though one new value (\"202\u00d73=\") duplicates an original\n# value used elsewhere in the table.\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"514\u00d76=\", \"985\u00d74=\"),\n    @(\"846\u00d75=\", \"299\u00d74=\"),\n    @(\"169\u00d78=\", \"369\u00d77=\"),\n    @(\"123\u00d76=\", \"937\u00d74=\"),\n    @(\"202\u00d73=\", \"831\u00d77=\"),\n    @(\"544\u00d73=\", \"515\u00d75=\"),\n    @(\"975\u00d77=\", \"402\u00d79=\"),\n    @(\"411\u00d72=\", \"321\u00d77=\"),\n    @(\"474\u00d78=\", \"965\u00d73=\"),\n    @(\"930\u00d73=\", \"868\u00d78=\"),\n    @(\"939\u00d74=\", \"916\u00d77=\"),\n    @(\"114\u00d78=\", \"718\u00d77=\"),\n    @(\"144\u00d78=\", \"110\u00d73=\"),\n    @(\"570\u00d76=\", \"437\u00d79=\"),\n    @(\"325\u00d74=\", \"818\u00d79=\"),\n    @(\"788\u00d78=\", \"721\u00d73=\"),\n    @(\"895\u00d76=\", \"310\u00d79=\"),\n    @(\"130\u00d73=\", \"606\u00d76=\"),\n    @(\"647\u00d75=\", \"466\u00d76=\"),\n    @(\"340\u00d79=\", \"202\u00d73=\"),\n    @(\"944\u00d78=\", \"575\u00d76=\"),\n    @(\"250\u00d76=\", \"555\u00d76=\"),\n    @(\"757\u00d72=\", \"668\u00d79=\"),\n    @(\"556\u00d77=\", \"679\u00d72=\"),\n    @(\"167\u00d79=\", \"576\u00d76=\")\n)\n\nforeach ($pair in $replacements) {\n    $fromText = $pair[0]\n    $toText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $fromText\n    $find.Replacement.Text = $toText\n    $find.Forward = $true\n    $find.Format = $false\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    $find.MatchSoundsLike = $false\n    $find.MatchAllWordForms = $false\n\n    # 0 = wdFindStop, 1 = wdReplaceOne: replace exactly the single\n    # (unique) occurrence of $fromText found in the document.\n    $find.Execute($fromText, $false, $false, $false, $false, $false, $true, 0, $false, $toText, 1) | Out-Null\n}\n"}
